# Regenerate save_data to use K (strike count) instead of Strike#
# Update column G ("K") values for rows 2-15 per recalculated std/mean s_vals.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0
$ws.Range("G3").Value = 2
$ws.Range("G4").Value = 3
$ws.Range("G5").Value = 2
$ws.Range("G6").Value = 1
$ws.Range("G7").Value = 0
$ws.Range("G8").Value = 0
$ws.Range("G9").Value = 1
$ws.Range("G10").Value = 0
$ws.Range("G12").Value = 2
$ws.Range("G13").Value = 2
$ws.Range("G15").Value = 2
